# Updated documentation and testing procedures (RTM workbook, Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 11 (test case #6) scenario/steps text was rewritten
$ws.Range("B11").Value = "Register Users after Log in"
$ws.Range("C11").Value = "1) Log In`r`n2) Log Out`r`n3) Register New User"

# The view scrolled down (row 10 now at the top) with C11 as the active cell
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C11").Select()
